$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1634.75
$ws.Range("I15").Value = 1634.75
$ws.Range("K15").Value = 4904.25
$ws.Range("M15").Value = -4735.25
$ws.Range("H19").Value = 775.1429000000001
$ws.Range("I19").Value = 943.6667
$ws.Range("K19").Value = 943.6667
$ws.Range("M19").Value = -768.6667
$ws.Range("H41").Value = 2483.125
$ws.Range("I41").Value = 2157.5
$ws.Range("K41").Value = 2157.5
$ws.Range("M41").Value = -1717.5
$ws.Range("H64").Value = 8574.583000000001
$ws.Range("J64").Value = 8832.777
$ws.Range("L64").Value = 8832.777
$ws.Range("N64").Value = -9328.777
$ws.Range("H67").Value = 8574.583000000001
$ws.Range("J67").Value = 8832.777
$ws.Range("L67").Value = 8832.777
$ws.Range("N67").Value = -10548.777
$ws.Range("H69").Value = 52401.184
$ws.Range("I69").Value = 89995
$ws.Range("J69").Value = 44047
$ws.Range("K69").Value = 269985
$ws.Range("L69").Value = 132141
$ws.Range("M69").Value = -269111
$ws.Range("N69").Value = -133889
$ws.Range("H72").Value = 52401.184
$ws.Range("I72").Value = 89995
$ws.Range("J72").Value = 44047
$ws.Range("K72").Value = 809955
$ws.Range("L72").Value = 396423
$ws.Range("M72").Value = -805587
$ws.Range("N72").Value = -405159
$ws.Range("H87").Value = 29999
$ws.Range("J87").Value = 29999
$ws.Range("L87").Value = 29999
$ws.Range("N87").Value = -32495
$ws.Range("H90").Value = 29999
$ws.Range("J90").Value = 29999
$ws.Range("L90").Value = 89997
$ws.Range("N90").Value = -102477
$ws.Range("H92").Value = 5436145.5
$ws.Range("I92").Value = 1225.8334
$ws.Range("J92").Value = 25001856
$ws.Range("K92").Value = 1225.8334
$ws.Range("L92").Value = 25001856
$ws.Range("M92").Value = 22.16660000000002
$ws.Range("N92").Value = -25004352
$ws.Range("H137").Value = 3415.6177
$ws.Range("I137").Value = 2245.4546
$ws.Range("J137").Value = 5560.9165
$ws.Range("K137").Value = 6736.3638
$ws.Range("L137").Value = 16682.7495
$ws.Range("M137").Value = -4186.3638
$ws.Range("N137").Value = -21782.7495
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2349.34
$ws.Range("I32").Value = 2349.34
$ws.Range("K32").Value = 2349.34
$ws.Range("M32").Value = -2062.34
$ws.Range("H101").Value = 42666.668
$ws.Range("J101").Value = 42666.668
$ws.Range("L101").Value = 42666.668
$ws.Range("N101").Value = -49156.668
$ws.Range("H112").Value = 36398.96
$ws.Range("J112").Value = 36398.96
$ws.Range("L112").Value = 36398.96
$ws.Range("N112").Value = -39352.96
$ws.Range("H122").Value = 3258.5454
$ws.Range("I122").Value = 2425.28
$ws.Range("K122").Value = 7275.84
$ws.Range("M122").Value = -4825.84
$ws.Range("H132").Value = 3878.1555
$ws.Range("I132").Value = 3458.7
$ws.Range("K132").Value = 10376.1
$ws.Range("M132").Value = -7846.099999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7145040
$ws.Range("I94").Value = 1753.3
$ws.Range("K94").Value = 1753.3
$ws.Range("M94").Value = -1302.3
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H134").Value = 4025.0938
$ws.Range("I134").Value = 3491.6191
$ws.Range("K134").Value = 10474.8573
$ws.Range("M134").Value = -7939.8573
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3454.18
$ws.Range("J31").Value = 4602.5713
$ws.Range("L31").Value = 4602.5713
$ws.Range("N31").Value = -5192.5713
$ws.Range("H34").Value = 3454.18
$ws.Range("J34").Value = 4602.5713
$ws.Range("L34").Value = 4602.5713
$ws.Range("N34").Value = -5006.5713
$ws.Range("H62").Value = 7266.6665
$ws.Range("I62").Value = 6424.75
$ws.Range("J62").Value = 8950.5
$ws.Range("K62").Value = 6424.75
$ws.Range("L62").Value = 8950.5
$ws.Range("M62").Value = -5800.75
$ws.Range("N62").Value = -10198.5
$ws.Range("H65").Value = 7266.6665
$ws.Range("I65").Value = 6424.75
$ws.Range("J65").Value = 8950.5
$ws.Range("K65").Value = 32123.75
$ws.Range("L65").Value = 44752.5
$ws.Range("M65").Value = -29003.75
$ws.Range("N65").Value = -50992.5
$ws.Range("H122").Value = 1253.7222
$ws.Range("I122").Value = 1175.75
$ws.Range("J122").Value = 1877.5
$ws.Range("K122").Value = 3527.25
$ws.Range("L122").Value = 5632.5
$ws.Range("M122").Value = -1077.25
$ws.Range("N122").Value = -10532.5
$ws.Range("H134").Value = 5349.1577
$ws.Range("I134").Value = 3103.5
$ws.Range("J134").Value = 7844.3335
$ws.Range("K134").Value = 9310.5
$ws.Range("L134").Value = 23533.0005
$ws.Range("M134").Value = -6775.5
$ws.Range("N134").Value = -28603.0005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 68468360
$ws.Range("I4").Value = 114111750
$ws.Range("J4").Value = 3262.1667
$ws.Range("K4").Value = 342335250
$ws.Range("L4").Value = 9786.500100000001
$ws.Range("M4").Value = -342335138
$ws.Range("N4").Value = -10010.5001
$ws.Range("H21").Value = 80.333336
$ws.Range("I21").Value = 95.5
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 286.5
$ws.Range("L21").Value = 150
$ws.Range("M21").Value = -113.5
$ws.Range("N21").Value = -496
$ws.Range("H24").Value = 3002
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 3002
$ws.Range("K24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("M24").Value = 9006
$ws.Range("N24").Value = -9466
$ws.Range("H26").Value = 191.44444
$ws.Range("I26").Value = 202.875
$ws.Range("K26").Value = 608.625
$ws.Range("M26").Value = -320.625
$ws.Range("H116").Value = 3532.4443
$ws.Range("J116").Value = 3532.4443
$ws.Range("L116").Value = 10597.3329
$ws.Range("N116").Value = -17481.3329
$ws.Range("H122").Value = 2086.3333
$ws.Range("J122").Value = 2230.1
$ws.Range("L122").Value = 20070.9
$ws.Range("N122").Value = -24970.9
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9467416
$ws.Range("I11").Value = 11336100
$ws.Range("K11").Value = 11336100
$ws.Range("M11").Value = -11335961
$ws.Range("I80").Value = 113805.6
$ws.Range("K80").Value = 113805.6
$ws.Range("M80").Value = -112807.6
$ws.Range("I83").Value = 113805.6
$ws.Range("K83").Value = 569028
$ws.Range("M83").Value = -564036
$ws.Range("H122").Value = 12587.608
$ws.Range("I122").Value = 17717.385
$ws.Range("J122").Value = 5918.9
$ws.Range("K122").Value = 53152.155
$ws.Range("L122").Value = 17756.7
$ws.Range("M122").Value = -50702.155
$ws.Range("N122").Value = -22656.7
$ws.Range("H126").Value = 3642.75
$ws.Range("I126").Value = 2472.3572
$ws.Range("K126").Value = 7417.071599999999
$ws.Range("M126").Value = -4947.071599999999
$ws.Range("H132").Value = 4659.9067
$ws.Range("I132").Value = 4760.5
$ws.Range("K132").Value = 14281.5
$ws.Range("M132").Value = -11751.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 2567.2856
$ws.Range("I43").Value = 1495.1666
$ws.Range("K43").Value = 1495.1666
$ws.Range("M43").Value = -1302.1666
$ws.Range("H82").Value = 1632.6
$ws.Range("I82").Value = 1413.3636
$ws.Range("K82").Value = 1413.3636
$ws.Range("M82").Value = -1052.3636
$ws.Range("H85").Value = 1632.6
$ws.Range("I85").Value = 1413.3636
$ws.Range("K85").Value = 1413.3636
$ws.Range("M85").Value = -165.3635999999999
$ws.Range("H100").Value = 48752.207
$ws.Range("I100").Value = 70491.44
$ws.Range("K100").Value = 70491.44
$ws.Range("M100").Value = -69950.44
$ws.Range("H110").Value = 100000
$ws.Range("I110").Value = 100000
$ws.Range("K110").Value = 100000
$ws.Range("M110").Value = -95910
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 20836814
$ws.Range("J81").Value = 83338090
$ws.Range("L81").Value = 166676180
$ws.Range("N81").Value = -166678302
$ws.Range("H84").Value = 20836814
$ws.Range("J84").Value = 83338090
$ws.Range("L84").Value = 833380900
$ws.Range("N84").Value = -833391508
$ws.Range("H113").Value = 502.48486
$ws.Range("I113").Value = 589.7368
$ws.Range("J113").Value = 384.07144
$ws.Range("K113").Value = 1769.2104
$ws.Range("L113").Value = 1152.21432
$ws.Range("M113").Value = 400.7896000000001
$ws.Range("N113").Value = -5492.21432
$ws.Range("H132").Value = 2585.375
$ws.Range("I132").Value = 1412.7693
$ws.Range("K132").Value = 4238.3079
$ws.Range("M132").Value = -1708.3079
